$wb = $excel.ActiveWorkbook
# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 11432
$ws.Range("F3").Value = 1955
$ws.Range("F4").Value = 532
$ws.Range("F5").Value = 854
$ws.Range("F6").Value = 2428
$ws.Range("F7").Value = 771
$ws.Range("F8").Value = 1017
$ws.Range("F9").Value = 607
$ws.Range("F10").Value = 465
$ws.Range("F11").Value = 1355
$ws.Range("F12").Value = 686
$ws.Range("F13").Value = 125
$ws.Range("F14").Value = 16
$ws.Range("F15").Value = 998
$ws.Range("F16").Value = 543
$ws.Range("F17").Value = 684
$ws.Range("F18").Value = 1135
$ws.Range("F19").Value = 215
$ws.Range("F20").Value = 942
$ws.Range("F21").Value = 12
$ws.Range("F22").Value = 140
$ws.Range("F23").Value = 312
$ws.Range("F25").Value = 267
$ws.Range("F28").Value = 685
$ws.Range("F29").Value = 181
$ws.Range("F30").Value = 116
$ws.Range("F31").Value = 328

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 904
$ws.Range("F7").Value = 9
$ws.Range("F8").Value = 105
$ws.Range("F10").Value = 399

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 59

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 11432
$ws.Range("F3").Value = 1955
$ws.Range("F5").Value = 532
$ws.Range("F6").Value = 854
$ws.Range("F7").Value = 2428
$ws.Range("F8").Value = 771
$ws.Range("F9").Value = 1017
$ws.Range("F11").Value = 607
$ws.Range("F12").Value = 465
$ws.Range("F13").Value = 59
$ws.Range("F14").Value = 1355
$ws.Range("F16").Value = 686
$ws.Range("F17").Value = 125
$ws.Range("F18").Value = 904
$ws.Range("F19").Value = 16
$ws.Range("F20").Value = 998
$ws.Range("F21").Value = 543
$ws.Range("F22").Value = 684
$ws.Range("F23").Value = 1135
$ws.Range("F24").Value = 215
$ws.Range("F25").Value = 942
$ws.Range("F26").Value = 12
$ws.Range("F27").Value = 140
$ws.Range("F28").Value = 312
$ws.Range("F31").Value = 267
$ws.Range("F32").Value = 9
$ws.Range("F33").Value = 105
$ws.Range("F34").Value = 105
$ws.Range("F37").Value = 685
$ws.Range("F38").Value = 181
$ws.Range("F40").Value = 116
$ws.Range("F41").Value = 399
$ws.Range("F43").Value = 328

